$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.042.65'
$ws.Range('E2').Value = '  +1.03%  '

$ws.Range('D3').Value = '1.568.48'
$ws.Range('E3').Value = '  +2.04%  '

$ws.Range('D4').Value = '''1.01'
$ws.Range('E4').Value = '  +0.35%  '

$ws.Range('D5').Value = '''207.71'
$ws.Range('E5').Value = '  +1.15%  '

$ws.Range('E6').Value = '  +1.31%  '

$ws.Range('D7').Value = '''1.01'
$ws.Range('E7').Value = '  +0.33%  '

$ws.Range('D8').Value = '''22.04'
$ws.Range('E8').Value = '  +3.89%  '

$ws.Range('D9').Value = '''0.250'
$ws.Range('E9').Value = '  +1.61%  '

$ws.Range('D10').Value = '''0.0588'
$ws.Range('E10').Value = '  +1.25%  '

$ws.Range('D11').Value = '''0.0859'
$ws.Range('E11').Value = '  +0.48%  '

$ws.Range('D12').Value = '1.789.71'
$ws.Range('E12').Value = '  +1.92%  '

$ws.Range('D13').Value = '1.589.78'
$ws.Range('E13').Value = '  +3.49%  '

$ws.Range('D14').Value = '''3.76'
$ws.Range('E14').Value = '  +2.74%  '

$ws.Range('D15').Value = '''0.522'
$ws.Range('E15').Value = '  +2.73%  '

$ws.Range('D16').Value = '27.049.57'
$ws.Range('E16').Value = '  +1.07%  '

$ws.Range('D17').Value = '''61.97'
$ws.Range('E17').Value = '  +1.79%  '

$ws.Range('D18').Value = '''218.49'
$ws.Range('E18').Value = '  +2.62%  '

$ws.Range('D19').Value = '0.0₃0697'
$ws.Range('E19').Value = '  +2.40%  '

$ws.Range('D20').Value = '''7.33'
$ws.Range('E20').Value = '  +1.40%  '

$ws.Range('E21').Value = '  +0.27%  '

$ws.Range('D22').Value = '''4.07'
$ws.Range('E22').Value = '  +1.71%  '

$ws.Range('E23').Value = '  +1.80%  '

$ws.Range('E24').Value = '  +1.17%  '

$ws.Range('D25').Value = '''154.40'
$ws.Range('E25').Value = '  +1.67%  '

$ws.Range('D26').Value = '''6.61'
$ws.Range('E26').Value = '  +0.70%  '

$ws.Range('D27').Value = '''14.98'
$ws.Range('E27').Value = '  +1.34%  '

$ws.Range('D28').Value = '''1.01'
$ws.Range('E28').Value = '  +0.35%  '

$ws.Range('E29').Value = '  +1.61%  '

$ws.Range('E30').Value = '  +3.32%  '

$ws.Range('E31').Value = '  +0.83%  '

$ws.Range('E32').Value = '  +0.49%  '

$ws.Range('D33').Value = '1.455.73'
$ws.Range('E33').Value = '  +6.79%  '

$ws.Range('D34').Value = '''3.06'
$ws.Range('E34').Value = '  +4.62%  '

$ws.Range('E35').Value = '  +4.46%  '

$ws.Range('D36').Value = '''0.966'
$ws.Range('E36').Value = '  +0.43%  '

$ws.Range('E37').Value = '  +0.74%  '

$ws.Range('E38').Value = '  +0.72%  '

$ws.Range('D39').Value = '''0.522'
$ws.Range('E39').Value = '  +0.46%  '

$ws.Range('D40').Value = '''0.815'
$ws.Range('E40').Value = '  +1.63%  '

$ws.Range('D41').Value = '''5.77'
$ws.Range('E41').Value = '  +0.39%  '

$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '''2.36'
$ws.Range('E42').Value = '  +7.42%  '

$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.31%  '

$ws.Range('D44').Value = '''0.988'
$ws.Range('E44').Value = '  -0.74%  '

$ws.Range('D45').Value = '''64.44'
$ws.Range('E45').Value = '  +2.51%  '

$ws.Range('D46').Value = '''1.77'
$ws.Range('E46').Value = '  +2.39%  '

$ws.Range('D47').Value = '1.705.51'
$ws.Range('E47').Value = '  +2.11%  '

$ws.Range('D48').Value = '''86.60'
$ws.Range('E48').Value = '  +1.73%  '

$ws.Range('D49').Value = '''0.0523'
$ws.Range('E49').Value = '  +2.95%  '

$ws.Range('E50').Value = '  +2.61%  '

$ws.Range('E51').Value = '  +0.52%  '
